$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label in N22 (literal text, including the trailing backtick, as committed)
$ws.Cells.Item(22, 14).Value = "w``"

# J64 keeps displaying "Number of ex" (re-entered by the author; same text)
$ws.Cells.Item(64, 10).Value = "Number of ex"

function Set-RowValues($row, $startCol, $values) {
    $col = $startCol
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col++
    }
}

# New cost-of-distributed-exchanges block, rows 101-105
Set-RowValues 101 2 @(16, 32, 64, 128, 256, 512, 1024, 2048, 4096)
Set-RowValues 102 2 @(4, 8, 16, 32, 64, 128, 256, 512, 1024)

$ws.Cells.Item(103, 1).Value = "sync"
Set-RowValues 103 2 @(32, 64, 128, 256, 512, 1024, 2048, 4096)
$ws.Cells.Item(103, 10).Formula = "=4096*2"

$ws.Cells.Item(104, 1).Value = "cent"
Set-RowValues 104 2 @(32, 64, 128, 256, 512, 1024, 2048, 4096)
$ws.Cells.Item(104, 10).Formula = "=4096*2"

$ws.Cells.Item(105, 1).Value = "decent"
Set-RowValues 105 2 @(0, 0, 0, 0, 0, 0, 0, 0, 0)

# Scroll / selection state mirrors the author's final view
$excel.ActiveWindow.ScrollRow = 68
$ws.Range("K105").Select()
